$d = $word.ActiveDocument

# Merge the split title runs ("Answers:" " " "Trigonometry" " " "(radians)")
# into a single run's text.
$d.Content.Find.Execute("Answers: Trigonometry (radians)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Answers: Trigonometry (radians)", 2)

# Merge the split abstract runs into a single run's text.
$d.Content.Find.Execute("Answers to the questions on trigonometry, using radians to measure angles.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Answers to the questions on trigonometry, using radians to measure angles.", 2)
